$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns, matching the style of the existing header row (AB1)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 54   # AC
    $ws.Cells.Item($r, 30).Value = 61   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}
